$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "sunny weather"
$ws.Cells.Item(2,2).Value = "晴れ|はれ"
$ws.Cells.Item(3,1).Value = "rain"
$ws.Cells.Item(3,2).Value = "雨|あめ"
$ws.Cells.Item(4,1).Value = "cloudy weather"
$ws.Cells.Item(4,2).Value = "曇り|くもり"
$ws.Cells.Item(5,1).Value = "snow"
$ws.Cells.Item(5,2).Value = "雪|ゆき"
$ws.Cells.Item(6,1).Value = "weather forecast"
$ws.Cells.Item(6,2).Value = "天気予報|てんきよほう"
$ws.Cells.Item(7,1).Value = "temperature (weather)"
$ws.Cells.Item(7,2).Value = "気温|きおん"
$ws.Cells.Item(8,1).Value = "summer"
$ws.Cells.Item(8,2).Value = "夏|なつ"
$ws.Cells.Item(9,1).Value = "winter"
$ws.Cells.Item(9,2).Value = "冬|ふゆ"
$ws.Cells.Item(10,1).Value = "this morning"
$ws.Cells.Item(10,2).Value = "今朝|けさ"
$ws.Cells.Item(11,1).Value = "the day after tomorrow"
$ws.Cells.Item(11,2).Value = "あさって"
$ws.Cells.Item(12,1).Value = "every week"
$ws.Cells.Item(12,2).Value = "毎週|まいしゅう"
$ws.Cells.Item(13,1).Value = "this month"
$ws.Cells.Item(13,2).Value = "今月|こんげつ"
$ws.Cells.Item(14,1).Value = "next month"
$ws.Cells.Item(14,2).Value = "来月|らいげつ"
$ws.Cells.Item(15,1).Value = "office worker"
$ws.Cells.Item(15,2).Value = "会社員|かいしゃいん"
$ws.Cells.Item(16,1).Value = "job; work; occupation"
$ws.Cells.Item(16,2).Value = "仕事|しごと"
$ws.Cells.Item(17,1).Value = "camera"
$ws.Cells.Item(17,2).Value = "カメラ"
$ws.Cells.Item(18,1).Value = "karaoke"
$ws.Cells.Item(18,2).Value = "カラオケ"
$ws.Cells.Item(19,1).Value = "place"
$ws.Cells.Item(19,2).Value = "所|ところ"
$ws.Cells.Item(20,1).Value = "tomato"
$ws.Cells.Item(20,2).Value = "トマト"
$ws.Cells.Item(21,1).Value = "chopsticks"
$ws.Cells.Item(21,2).Value = "はし"
$ws.Cells.Item(22,1).Value = "party"
$ws.Cells.Item(22,2).Value = "パーティー"
$ws.Cells.Item(23,1).Value = "barbecue"
$ws.Cells.Item(23,2).Value = "バーベキュー"
$ws.Cells.Item(24,1).Value = "homestay; living with a local family"
$ws.Cells.Item(24,2).Value = "ホームステイ"
$ws.Cells.Item(25,1).Value = "bath"
$ws.Cells.Item(25,2).Value = "お風呂|おふろ"
$ws.Cells.Item(26,1).Value = "Spain"
$ws.Cells.Item(26,2).Value = "スペイン"
$ws.Cells.Item(27,1).Value = "something"
$ws.Cells.Item(27,2).Value = "何か|なにか"
$ws.Cells.Item(28,1).Value = "skillful; good at (～が)"
$ws.Cells.Item(28,2).Value = "上手|じょうず(な)"
$ws.Cells.Item(29,1).Value = "clumsy; poor at (～が)"
$ws.Cells.Item(29,2).Value = "下手|へた(な)"
$ws.Cells.Item(30,1).Value = "famous"
$ws.Cells.Item(30,2).Value = "有名|ゆうめい(な)"
$ws.Cells.Item(31,1).Value = "to wash (～を)"
$ws.Cells.Item(31,2).Value = "洗う|あらう"
$ws.Cells.Item(32,1).Value = "to say"
$ws.Cells.Item(32,2).Value = "言う|いう"
$ws.Cells.Item(33,1).Value = "to need (～が)"
$ws.Cells.Item(33,2).Value = "いる"
$ws.Cells.Item(34,1).Value = "to be late"
$ws.Cells.Item(34,2).Value = "遅くなる|おそくなる"
$ws.Cells.Item(35,1).Value = "to take a bath"
$ws.Cells.Item(35,2).Value = "お風呂に入る|おふろにはいる"
$ws.Cells.Item(36,1).Value = "to think"
$ws.Cells.Item(36,2).Value = "思う|おもう"
$ws.Cells.Item(37,1).Value = "to cut (～を)"
$ws.Cells.Item(37,2).Value = "切る|きる"
$ws.Cells.Item(38,1).Value = "to make (～を)"
$ws.Cells.Item(38,2).Value = "作る|つくる"
$ws.Cells.Item(39,1).Value = "rain falls"
$ws.Cells.Item(39,2).Value = "雨が降る|あめがふる"
$ws.Cells.Item(40,1).Value = "snow falls"
$ws.Cells.Item(40,2).Value = "雪が降る|ゆきがふる"
$ws.Cells.Item(41,1).Value = "to take (a thing) (～を)"
$ws.Cells.Item(41,2).Value = "持っていく|もっていく"
$ws.Cells.Item(42,1).Value = "to throw away (～を)"
$ws.Cells.Item(42,2).Value = "捨てる|すてる"
$ws.Cells.Item(43,1).Value = "to begin (～を)"
$ws.Cells.Item(43,2).Value = "始める|はじめる"
$ws.Cells.Item(44,1).Value = "to drive (～を)"
$ws.Cells.Item(44,2).Value = "運転する|うんてんする"
$ws.Cells.Item(45,1).Value = "to do laundry (～を)"
$ws.Cells.Item(45,2).Value = "洗濯する|せんたくする"
$ws.Cells.Item(46,1).Value = "to clean (～を)"
$ws.Cells.Item(46,2).Value = "掃除する|そうじする"
$ws.Cells.Item(47,1).Value = "to cook"
$ws.Cells.Item(47,2).Value = "料理する|りょうりする"
$ws.Cells.Item(48,1).Value = "uh-huh; yes"
$ws.Cells.Item(48,2).Value = "うん"
$ws.Cells.Item(49,1).Value = "uh-uh; no"
$ws.Cells.Item(49,2).Value = "ううん"
$ws.Cells.Item(50,1).Value = "always"
$ws.Cells.Item(50,2).Value = "いつも"
$ws.Cells.Item(51,1).Value = "(do something) late"
$ws.Cells.Item(51,2).Value = "遅く|おそく"
$ws.Cells.Item(52,1).Value = "Cheers! (a toast)"
$ws.Cells.Item(52,2).Value = "乾杯|かんぱい"
$ws.Cells.Item(53,1).Value = "all (of the people) together"
$ws.Cells.Item(53,2).Value = "みんなで"
$ws.Cells.Item(54,1).Value = "That's too bad."
$ws.Cells.Item(54,2).Value = "残念(ですね)|ざんねん(ですね)"
$ws.Cells.Item(55,1).Value = "not...yet (w/negative)"
$ws.Cells.Item(55,2).Value = "まだ"
$ws.Cells.Item(56,1).Value = "about...; concerning..."
$ws.Cells.Item(56,2).Value = "～について"
$ws.Cells.Item(57,1).Value = "...degrees (temperature)"
$ws.Cells.Item(57,2).Value = "～度|～ど"
$ws.Cells.Item(58,1).Value = "how"
$ws.Cells.Item(58,2).Value = "どう"
$ws.Cells.Item(59,1).Value = "rice"
$ws.Cells.Item(59,2).Value = "ご飯|ごはん"
$ws.Cells.Item(60,1).Value = "side dish"
$ws.Cells.Item(60,2).Value = "おかず"
$ws.Cells.Item(61,1).Value = "miso soup"
$ws.Cells.Item(61,2).Value = "みそ汁|みそしる"
$ws.Cells.Item(62,1).Value = "set meal"
$ws.Cells.Item(62,2).Value = "定食|ていしょく"
$ws.Cells.Item(63,1).Value = "curry with rice"
$ws.Cells.Item(63,2).Value = "カレーライス"
$ws.Cells.Item(64,1).Value = "rice balls"
$ws.Cells.Item(64,2).Value = "おにぎり"
$ws.Cells.Item(65,1).Value = "ramen noodles"
$ws.Cells.Item(65,2).Value = "ラーメン"
$ws.Cells.Item(66,1).Value = "udon noodles"
$ws.Cells.Item(66,2).Value = "うどん"
$ws.Cells.Item(67,1).Value = "pasta"
$ws.Cells.Item(67,2).Value = "パスタ"
$ws.Cells.Item(68,1).Value = "dumplings"
$ws.Cells.Item(68,2).Value = "ぎょうざ"
$ws.Cells.Item(69,1).Value = "beef rice bowl"
$ws.Cells.Item(69,2).Value = "牛丼|ぎゅうどん"
$ws.Cells.Item(70,1).Value = "hamburger steak"
$ws.Cells.Item(70,2).Value = "ハンバーグ"
$ws.Cells.Item(71,1).Value = "raw seafood"
$ws.Cells.Item(71,2).Value = "さしみ"
$ws.Cells.Item(72,1).Value = "savory pancake"
$ws.Cells.Item(72,2).Value = "お好み焼き|おこのみやき"
$ws.Cells.Item(73,1).Value = "toast"
$ws.Cells.Item(73,2).Value = "トースト"
$ws.Cells.Item(74,1).Value = "soup"
$ws.Cells.Item(74,2).Value = "スープ"
$ws.Cells.Item(75,1).Value = "yogurt"
$ws.Cells.Item(75,2).Value = "ヨーグルト"
$ws.Cells.Item(76,1).Value = "broiled fish"
$ws.Cells.Item(76,2).Value = "焼き魚|やきざかな"
$ws.Cells.Item(77,1).Value = "egg"
$ws.Cells.Item(77,2).Value = "たまご"
$ws.Cells.Item(78,1).Value = "Do you have avocado?"
$ws.Cells.Item(78,2).Value = "アボカドはありますか。"
$ws.Cells.Item(79,1).Value = "Is there alcohol in this?"
$ws.Cells.Item(79,2).Value = "この中にお酒が入っていますか。|このなかにおさけがはいっていますか。"
$ws.Cells.Item(80,1).Value = "Is this halal?"
$ws.Cells.Item(80,2).Value = "これはハラルフードですか。"
$ws.Cells.Item(81,1).Value = "I have an allergy to peanuts."
$ws.Cells.Item(81,2).Value = "ピーナッツアレルギーがあります。"
$ws.Cells.Item(82,1).Value = "Can I have a bag?"
$ws.Cells.Item(82,2).Value = "袋をお願いします。|ふくろをおねがいします。"
$ws.Cells.Item(83,1).Value = "fruit"
$ws.Cells.Item(83,2).Value = "果物|くだもの"
$ws.Cells.Item(84,1).Value = "strawberry"
$ws.Cells.Item(84,2).Value = "いちご"
$ws.Cells.Item(85,1).Value = "watermelon"
$ws.Cells.Item(85,2).Value = "すいか"
$ws.Cells.Item(86,1).Value = "mandarin orange"
$ws.Cells.Item(86,2).Value = "みかん"
$ws.Cells.Item(87,1).Value = "apple"
$ws.Cells.Item(87,2).Value = "りんご"
$ws.Cells.Item(88,1).Value = "peach"
$ws.Cells.Item(88,2).Value = "もも"
$ws.Cells.Item(89,1).Value = "grapes"
$ws.Cells.Item(89,2).Value = "ぶどう"
$ws.Cells.Item(90,1).Value = "vegetable"
$ws.Cells.Item(90,2).Value = "野菜|やさい"
$ws.Cells.Item(91,1).Value = "carrot"
$ws.Cells.Item(91,2).Value = "にんじん"
$ws.Cells.Item(92,1).Value = "onion"
$ws.Cells.Item(92,2).Value = "たまねぎ"
$ws.Cells.Item(93,1).Value = "potato"
$ws.Cells.Item(93,2).Value = "じゃがいも"
$ws.Cells.Item(94,1).Value = "eggplant"
$ws.Cells.Item(94,2).Value = "なす"
$ws.Cells.Item(95,1).Value = "cucumber"
$ws.Cells.Item(95,2).Value = "きゅうり"
$ws.Cells.Item(96,1).Value = "cabbage"
$ws.Cells.Item(96,2).Value = "キャベツ"
$ws.Cells.Item(97,1).Value = "meat"
$ws.Cells.Item(97,2).Value = "肉|にく"
$ws.Cells.Item(98,1).Value = "beef"
$ws.Cells.Item(98,2).Value = "牛肉|ぎゅうにく"
$ws.Cells.Item(99,1).Value = "pork"
$ws.Cells.Item(99,2).Value = "豚肉|ぶたにく"
$ws.Cells.Item(100,1).Value = "chicken"
$ws.Cells.Item(100,2).Value = "鶏肉|とりにく"
$ws.Cells.Item(101,1).Value = "office worker"
$ws.Cells.Item(101,2).Value = "会社員|かいしゃいん"
$ws.Cells.Item(102,1).Value = "store clerk"
$ws.Cells.Item(102,2).Value = "店員|てんいん"
$ws.Cells.Item(103,1).Value = "member"
$ws.Cells.Item(103,2).Value = "会員|かいいん"
$ws.Cells.Item(104,1).Value = "station attendant"
$ws.Cells.Item(104,2).Value = "駅員|えきいん"
$ws.Cells.Item(105,1).Value = "new"
$ws.Cells.Item(105,2).Value = "新しい|あたらしい"
$ws.Cells.Item(106,1).Value = "newspaper"
$ws.Cells.Item(106,2).Value = "新聞|しんぶん"
$ws.Cells.Item(107,1).Value = "Bullet Train"
$ws.Cells.Item(107,2).Value = "新幹線|しんかんせん"
$ws.Cells.Item(108,1).Value = "fresh"
$ws.Cells.Item(108,2).Value = "新鮮な|しんせんな"
$ws.Cells.Item(109,1).Value = "to listen"
$ws.Cells.Item(109,2).Value = "聞く|きく"
$ws.Cells.Item(110,1).Value = "can be heard"
$ws.Cells.Item(110,2).Value = "聞こえる|きこえる"
$ws.Cells.Item(111,1).Value = "to make"
$ws.Cells.Item(111,2).Value = "作る|つくる"
$ws.Cells.Item(112,1).Value = "composition"
$ws.Cells.Item(112,2).Value = "作文|さくぶん"
$ws.Cells.Item(113,1).Value = "artistic piece"
$ws.Cells.Item(113,2).Value = "作品|さくひん"
$ws.Cells.Item(114,1).Value = "author"
$ws.Cells.Item(114,2).Value = "作者|さくしゃ"
$ws.Cells.Item(115,1).Value = "job"
$ws.Cells.Item(115,2).Value = "仕事|しごと"
$ws.Cells.Item(116,1).Value = "revenge"
$ws.Cells.Item(116,2).Value = "仕返し|しかえし"
$ws.Cells.Item(117,1).Value = "to serve; to work under"
$ws.Cells.Item(117,2).Value = "仕える|つかえる"
$ws.Cells.Item(118,1).Value = "thing"
$ws.Cells.Item(118,2).Value = "事|こと"
$ws.Cells.Item(119,1).Value = "fire"
$ws.Cells.Item(119,2).Value = "火事|かじ"
$ws.Cells.Item(120,1).Value = "meal"
$ws.Cells.Item(120,2).Value = "食事|しょくじ"
$ws.Cells.Item(121,1).Value = "reply"
$ws.Cells.Item(121,2).Value = "返事|へんじ"
$ws.Cells.Item(122,1).Value = "train"
$ws.Cells.Item(122,2).Value = "電車|でんしゃ"
$ws.Cells.Item(123,1).Value = "electricity; light"
$ws.Cells.Item(123,2).Value = "電気|でんき"
$ws.Cells.Item(124,1).Value = "telephone"
$ws.Cells.Item(124,2).Value = "電話|でんわ"
$ws.Cells.Item(125,1).Value = "battery"
$ws.Cells.Item(125,2).Value = "電池|でんち"
$ws.Cells.Item(126,1).Value = "microwave oven"
$ws.Cells.Item(126,2).Value = "電子レンジ|でんしレンジ"
$ws.Cells.Item(127,1).Value = "car"
$ws.Cells.Item(127,2).Value = "車|くるま"
$ws.Cells.Item(128,1).Value = "bicycle"
$ws.Cells.Item(128,2).Value = "自転車|じてんしゃ"
$ws.Cells.Item(129,1).Value = "wheelchair"
$ws.Cells.Item(129,2).Value = "車いす|くるまいす"
$ws.Cells.Item(130,1).Value = "parking lot"
$ws.Cells.Item(130,2).Value = "駐車場|ちゅうしゃじょう"
$ws.Cells.Item(131,1).Value = "to be absent; to rest"
$ws.Cells.Item(131,2).Value = "休む|やすむ"
$ws.Cells.Item(132,1).Value = "holiday; absence"
$ws.Cells.Item(132,2).Value = "休み|やすみ"
$ws.Cells.Item(133,1).Value = "holiday"
$ws.Cells.Item(133,2).Value = "休日|きゅうじつ"
$ws.Cells.Item(134,1).Value = "to say"
$ws.Cells.Item(134,2).Value = "言う|いう"
$ws.Cells.Item(135,1).Value = "linguistics"
$ws.Cells.Item(135,2).Value = "言語学|げんごがく"
$ws.Cells.Item(136,1).Value = "dialect"
$ws.Cells.Item(136,2).Value = "方言|ほうげん"
$ws.Cells.Item(137,1).Value = "word; language"
$ws.Cells.Item(137,2).Value = "言葉|ことば"
$ws.Cells.Item(138,1).Value = "to read"
$ws.Cells.Item(138,2).Value = "読む|よむ"
$ws.Cells.Item(139,1).Value = "reading books"
$ws.Cells.Item(139,2).Value = "読書|どくしょ"
$ws.Cells.Item(140,1).Value = "reading matter"
$ws.Cells.Item(140,2).Value = "読み物|よみもの"
$ws.Cells.Item(141,1).Value = "to think"
$ws.Cells.Item(141,2).Value = "思う|おもう"
$ws.Cells.Item(142,1).Value = "mysterious"
$ws.Cells.Item(142,2).Value = "不思議な|ふしぎな"
$ws.Cells.Item(143,1).Value = "to recall; to remember"
$ws.Cells.Item(143,2).Value = "思い出す|おもいだす"
$ws.Cells.Item(144,1).Value = "next"
$ws.Cells.Item(144,2).Value = "次|つぎ"
$ws.Cells.Item(145,1).Value = "second daughter"
$ws.Cells.Item(145,2).Value = "次女|じじょ"
$ws.Cells.Item(146,1).Value = "table of contents"
$ws.Cells.Item(146,2).Value = "目次|もくじ"
$ws.Cells.Item(147,1).Value = "next time"
$ws.Cells.Item(147,2).Value = "次回|じかい"
$ws.Cells.Item(148,1).Value = "what"
$ws.Cells.Item(148,2).Value = "何|なに"
$ws.Cells.Item(149,1).Value = "what time"
$ws.Cells.Item(149,2).Value = "何時|なんじ"
$ws.Cells.Item(150,1).Value = "how many people"
$ws.Cells.Item(150,2).Value = "何人|なんにん"
$ws.Cells.Item(151,1).Value = "something"
$ws.Cells.Item(151,2).Value = "何か|なにか"
